# Apply the "Add budget outputs with UD penalty" edit.
# Touches five sheets: Summary, Costs and Revenues, Fed-in Capacity,
# Unmet Demand, Household Surplus.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Summary sheet
# ---------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 0.01
$wsSummary.Range("B6").Value = -218750.5042824882
$wsSummary.Range("B7").Value = 5413711.842050619
$wsSummary.Range("B8").Value = 22726010.95505212
$wsSummary.Range("B10").Value = 4350924.71095404

# ---------------------------------------------------------------
# 2. Costs and Revenues sheet - row 6 ("Total Profits")
# ---------------------------------------------------------------
$wsCosts = $wb.Worksheets.Item("Costs and Revenues")
$wsCosts.Range("B6").Value = -48778.27397003476
$wsCosts.Range("C6").Value = -48778.27397003476
$wsCosts.Range("D6").Value = -48778.27397003476
$wsCosts.Range("E6").Value = -15150.67397003476
$wsCosts.Range("F6").Value = -15150.67397003476
$wsCosts.Range("G6").Value = -15150.67397003476
$wsCosts.Range("H6").Value = -15150.67397003476
$wsCosts.Range("I6").Value = -15150.67397003476
$wsCosts.Range("J6").Value = -15150.67397003476
$wsCosts.Range("K6").Value = -15150.67397003476
$wsCosts.Range("L6").Value = -15150.67397003476
$wsCosts.Range("M6").Value = -15150.67397003476
$wsCosts.Range("N6").Value = -15150.67397003476
$wsCosts.Range("O6").Value = -15150.67397003476
$wsCosts.Range("P6").Value = -15150.67397003476

# ---------------------------------------------------------------
# 3. Household Surplus sheet - column B, rows 2-16 all become the
#    same constant value
# ---------------------------------------------------------------
$wsHH = $wb.Worksheets.Item("Household Surplus")
for ($r = 2; $r -le 16; $r++) {
    $wsHH.Cells.Item($r, 2).Value = 376275.9193600624
}

# ---------------------------------------------------------------
# 4. Fed-in Capacity and Unmet Demand sheets
#
# Both sheets have a 45-row grid (rows 2-46) split into three
# repeating "day type" groups based on (row-2) mod 3. Within each
# group, columns I-R hold a canonical pattern; this edit fills in
# (Fed-in Capacity) or overwrites (Unmet Demand) those columns for
# every row in every group with the group's canonical values.
# ---------------------------------------------------------------

$g0_fedin = @{10=169.0966151720738; 11=220.0898510449805; 12=235.7664149699872; 13=230.3462332272727; 14=229.4130635965909; 15=230.0982114216867; 16=231.2329957552695; 17=212.3149906599047; 18=65.71641987298243}
$g1_fedin = @{9=10.12574714858493; 10=126.0910353404088; 11=137.841438974359; 12=138.5543797798742; 13=142.1340339220183; 14=131.3417120833333; 15=142.5962444444444; 16=133.9744074143302; 17=139.9817740860215; 18=45.52166981132082}
$g2_fedin = @{10=33.63624132272333; 11=106.7437663446525; 12=134.8846762812383; 13=138.9257839476051; 14=127.6855444652332; 15=138.4565384518428; 16=135.0065633140411; 17=65.34295837775146}

$g0_unmet = @{10=11.94928935461252; 11=0; 12=0; 13=0; 14=0; 15=0; 16=0; 17=9.990699214544804; 18=149.8691179411497}
$g1_unmet = @{9=89.39663285141508; 10=0.7465913262578567; 11=0; 12=0; 13=0; 14=0; 15=0; 16=0; 17=0; 18=100.1578341526431}
$g2_unmet = @{10=93.35918011667277; 11=22.26949182588285; 12=0; 13=0; 14=0; 15=0; 16=2.721440735106512; 17=86.16204325169439}

$wsFedin = $wb.Worksheets.Item("Fed-in Capacity")
$wsUnmet = $wb.Worksheets.Item("Unmet Demand")

for ($r = 2; $r -le 46; $r++) {
    $grp = ($r - 2) % 3
    if ($grp -eq 0) {
        $mapFedin = $g0_fedin
        $mapUnmet = $g0_unmet
    } elseif ($grp -eq 1) {
        $mapFedin = $g1_fedin
        $mapUnmet = $g1_unmet
    } else {
        $mapFedin = $g2_fedin
        $mapUnmet = $g2_unmet
    }

    foreach ($col in $mapFedin.Keys) {
        $wsFedin.Cells.Item($r, $col).Value = $mapFedin[$col]
    }
    foreach ($col in $mapUnmet.Keys) {
        $wsUnmet.Cells.Item($r, $col).Value = $mapUnmet[$col]
    }
}

Write-Host "Edit applied."
